# BagUserItemData: add a "quality" (itemNew) column to the item base csv sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- new column H: "quality" -----------------------------------------
# Header cell (row 1) gets a distinct "new field" font so it stands out
# from the rest of the header row.
$ws.Range("H1").Value = "quality"
$ws.Range("H1").Font.Size = 10
$ws.Range("H1").Font.Color = 11171480   # RGB(152,118,170) -> 0xAA7698 (BGR) == FF9876AA (ARGB)
$ws.Range("H1").Font.Name = "Arial Unicode MS"

# Row 2 = field type, row 3 = field description, row 4 = sample data,
# matching the existing A:G columns. These stay in the sheet's default
# (unstyled) format, same as how a plain value paste behaves.
$ws.Range("H2").Value = "int"
$ws.Range("H3").Value = "des"
$ws.Range("H4").Value = 1

# Leave the cursor on the last edited cell, like the author would have.
$ws.Range("H4").Select()
